$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 229.8
$ws.Range("I92").Value = 199.76471
$ws.Range("K92").Value = 199.76471
$ws.Range("M92").Value = 1048.23529
$ws.Range("H96").Value = 297.31818
$ws.Range("I96").Value = 246.89473
$ws.Range("J96").Value = 616.6667
$ws.Range("K96").Value = 740.6841900000001
$ws.Range("L96").Value = 1850.0001
$ws.Range("M96").Value = 632.3158099999999
$ws.Range("N96").Value = -4596.0001
$ws.Range("H98").Value = 1143.4348
$ws.Range("I98").Value = 1140.8636
$ws.Range("J98").Value = 1200
$ws.Range("K98").Value = 1140.8636
$ws.Range("L98").Value = 1200
$ws.Range("M98").Value = 357.1364000000001
$ws.Range("N98").Value = -4196
$ws.Range("H122").Value = 1143.4348
$ws.Range("I122").Value = 1140.8636
$ws.Range("J122").Value = 1200
$ws.Range("K122").Value = 3422.5908
$ws.Range("L122").Value = 3600
$ws.Range("M122").Value = -972.5907999999999
$ws.Range("N122").Value = -8500
$ws.Range("H132").Value = 2318.111
$ws.Range("I132").Value = 1724.7937
$ws.Range("J132").Value = 4394.722
$ws.Range("K132").Value = 5174.3811
$ws.Range("L132").Value = 13184.166
$ws.Range("M132").Value = -2644.3811
$ws.Range("N132").Value = -18244.166
$ws.Range("H135").Value = 205630.33
$ws.Range("I135").Value = 239548.33
$ws.Range("J135").Value = 2122.2856
$ws.Range("K135").Value = 2155934.97
$ws.Range("L135").Value = 19100.5704
$ws.Range("M135").Value = -2153399.97
$ws.Range("N135").Value = -24170.5704
$ws.Range("H137").Value = 2660.976
$ws.Range("I137").Value = 921.63635
$ws.Range("K137").Value = 2764.90905
$ws.Range("M137").Value = -214.9090500000002
$ws.Range("H138").Value = 2262.1226
$ws.Range("J138").Value = 3094.1206
$ws.Range("L138").Value = 9282.361800000001
$ws.Range("N138").Value = -19562.3618
$ws.Range("H141").Value = 2484.3804
$ws.Range("I141").Value = 2468.1875
$ws.Range("J141").Value = 2489.0908
$ws.Range("K141").Value = 7404.5625
$ws.Range("L141").Value = 7467.2724
$ws.Range("M141").Value = -2224.5625
$ws.Range("N141").Value = -17827.2724

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 35354.465
$ws.Range("I32").Value = 13103.305
$ws.Range("K32").Value = 13103.305
$ws.Range("M32").Value = -12816.305
$ws.Range("H61").Value = 1922.069
$ws.Range("I61").Value = 1779.2593
$ws.Range("K61").Value = 1779.2593
$ws.Range("M61").Value = -1567.2593
$ws.Range("H74").Value = 19394.281
$ws.Range("I74").Value = 1215.826
$ws.Range("K74").Value = 1215.826
$ws.Range("M74").Value = -341.826
$ws.Range("H77").Value = 19394.281
$ws.Range("I77").Value = 1215.826
$ws.Range("K77").Value = 6079.13
$ws.Range("M77").Value = -1711.13
$ws.Range("H136").Value = 1922.069
$ws.Range("I136").Value = 1779.2593
$ws.Range("K136").Value = 5337.7779
$ws.Range("M136").Value = -2787.7779

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H122").Value = 42465
$ws.Range("J122").Value = 42465
$ws.Range("L122").Value = 42465
$ws.Range("N122").Value = -52265

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3872.074
$ws.Range("I31").Value = 1946.9656
$ws.Range("J31").Value = 6105.2
$ws.Range("K31").Value = 1946.9656
$ws.Range("L31").Value = 6105.2
$ws.Range("M31").Value = -1651.9656
$ws.Range("N31").Value = -6695.2
$ws.Range("H34").Value = 3872.074
$ws.Range("I34").Value = 1946.9656
$ws.Range("J34").Value = 6105.2
$ws.Range("K34").Value = 1946.9656
$ws.Range("L34").Value = 6105.2
$ws.Range("M34").Value = -1744.9656
$ws.Range("N34").Value = -6509.2
$ws.Range("H58").Value = 1613.6296
$ws.Range("I58").Value = 1497.8096
$ws.Range("J58").Value = 2019
$ws.Range("K58").Value = 1497.8096
$ws.Range("L58").Value = 2019
$ws.Range("M58").Value = -1294.8096
$ws.Range("N58").Value = -2425
$ws.Range("H136").Value = 1613.6296
$ws.Range("I136").Value = 1497.8096
$ws.Range("J136").Value = 2019
$ws.Range("K136").Value = 4493.4288
$ws.Range("L136").Value = 6057
$ws.Range("M136").Value = -1943.4288
$ws.Range("N136").Value = -11157

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1133.3334
$ws.Range("I68").Value = 1000
$ws.Range("J68").Value = 1200
$ws.Range("K68").Value = 3000
$ws.Range("L68").Value = 3600
$ws.Range("M68").Value = -2189
$ws.Range("N68").Value = -5222
$ws.Range("H71").Value = 1133.3334
$ws.Range("I71").Value = 1000
$ws.Range("J71").Value = 1200
$ws.Range("K71").Value = 9000
$ws.Range("L71").Value = 10800
$ws.Range("M71").Value = -4944
$ws.Range("N71").Value = -18912
$ws.Range("H131").Value = 919.62103
$ws.Range("I131").Value = 457.5
$ws.Range("J131").Value = 962.1149
$ws.Range("K131").Value = 1372.5
$ws.Range("L131").Value = 2886.3447
$ws.Range("M131").Value = 3667.5
$ws.Range("N131").Value = -12966.3447

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 13068.5
$ws.Range("J46").Value = 16090.25
$ws.Range("L46").Value = 16090.25
$ws.Range("N46").Value = -16402.25
$ws.Range("H80").Value = 3129.4285
$ws.Range("I80").Value = 3000
$ws.Range("K80").Value = 3000
$ws.Range("M80").Value = -2002
$ws.Range("H83").Value = 3129.4285
$ws.Range("I83").Value = 3000
$ws.Range("K83").Value = 15000
$ws.Range("M83").Value = -10008
$ws.Range("H128").Value = 45690
$ws.Range("J128").Value = 45690
$ws.Range("L128").Value = 45690
$ws.Range("N128").Value = -55650
$ws.Range("H132").Value = 2553.849
$ws.Range("I132").Value = 2340.2104
$ws.Range("J132").Value = 3095.0667
$ws.Range("K132").Value = 7020.6312
$ws.Range("L132").Value = 9285.2001
$ws.Range("M132").Value = -4490.6312
$ws.Range("N132").Value = -14345.2001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H63").Value = 34085
$ws.Range("J63").Value = 34085
$ws.Range("L63").Value = 34085
$ws.Range("N63").Value = -35583
$ws.Range("H66").Value = 34085
$ws.Range("J66").Value = 34085
$ws.Range("L66").Value = 102255
$ws.Range("N66").Value = -109743
$ws.Range("H136").Value = 4256
$ws.Range("I136").Value = 1526.3103
$ws.Range("J136").Value = 13051.667
$ws.Range("K136").Value = 4578.9309
$ws.Range("L136").Value = 39155.001
$ws.Range("M136").Value = -2028.9309
$ws.Range("N136").Value = -44255.001
$ws.Range("H137").Value = 36695.8
$ws.Range("J137").Value = 40869.75
$ws.Range("L137").Value = 40869.75
$ws.Range("N137").Value = -51069.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 25000
$ws.Range("J92").Value = 25000
$ws.Range("L92").Value = 25000
$ws.Range("N92").Value = -29992
$ws.Range("H96").Value = 1533.3334
$ws.Range("I96").Value = 800
$ws.Range("J96").Value = 3000
$ws.Range("K96").Value = 800
$ws.Range("L96").Value = 3000
$ws.Range("M96").Value = 573
$ws.Range("N96").Value = -5746
$ws.Range("H109").Value = 17338.5
$ws.Range("J109").Value = 17338.5
$ws.Range("L109").Value = 17338.5
$ws.Range("N109").Value = -20112.5
$ws.Range("H122").Value = 8754.406000000001
$ws.Range("I122").Value = 10827.318
$ws.Range("J122").Value = 4194
$ws.Range("K122").Value = 32481.954
$ws.Range("L122").Value = 12582
$ws.Range("M122").Value = -30031.954
$ws.Range("N122").Value = -17482
